$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 7 (existing employee "PEDROZO DANTE OSCALDO") -----------------
# Dependencia de Revista changes from "art17 trab temp" to "PEON GRAL"
$ws.Range("F7").Value = "PEON GRAL"
# Situacion de Revista 1 flag goes from 1 to 0
$ws.Range("J7").Value = 0

# --- Row 8 (new employee "BRIZUELA JUAN ANDRES RAMON") ------------------
$ws.Range("B8").Value = "BRIZUELA JUAN ANDRES RAMON"
$ws.Range("C8").Value = 20412270089
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 32
$ws.Range("F8").Value = "PEON GRAL"
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 5
$ws.Range("Q8").Value = 97
$ws.Range("R8").Value = 111
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 27
$ws.Range("U8").Value = 1
$ws.Range("V8").Value = 1
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = 0
$ws.Range("Y8").Value = 0
$ws.Range("Z8").Value = 0
$ws.Range("AA8").Value = 30
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 0
$ws.Range("AD8").Value = 2
$ws.Range("AE8").Value = 119302
$ws.Range("AF8").Value = 0
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").Value = 0
$ws.Range("AJ8").Value = 0
$ws.Range("AK8").Value = 0
$ws.Range("AL8").Value = 0
$ws.Rows.Item(8).RowHeight = 30

# --- Row 9 (new employee "NAVARRO JANATAN FACUNDO") ----------------------
$ws.Range("B9").Value = "NAVARRO JANATAN FACUNDO"
$ws.Range("C9").Value = 20363899952
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 51
$ws.Range("F9").Value = "PEON GRAL"
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 5
$ws.Range("Q9").Value = 97
$ws.Range("R9").Value = 111
$ws.Range("S9").Value = 0
$ws.Range("T9").Value = 27
$ws.Range("U9").Value = 1
$ws.Range("V9").Value = 1
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = 0
$ws.Range("Y9").Value = 0
$ws.Range("Z9").Value = 0
$ws.Range("AA9").Value = 30
$ws.Range("AB9").Value = 0
$ws.Range("AC9").Value = 0
$ws.Range("AD9").Value = 2
$ws.Range("AE9").Value = 119302
$ws.Range("AF9").Value = 0
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0
$ws.Range("AI9").Value = 0
$ws.Range("AJ9").Value = 0
$ws.Range("AK9").Value = 0
$ws.Range("AL9").Value = 0
$ws.Rows.Item(9).RowHeight = 30

# --- Rows 10-19 (already-present blank employee rows): fill the tope/LRT
#     columns that are now populated with default values ------------------
for ($r = 10; $r -le 19; $r++) {
    $ws.Range("AD$r").Value = 2
    $ws.Range("AE$r").Value = 119302
    $ws.Range("AF$r").Value = 0
    $ws.Range("AG$r").Value = 0
    $ws.Range("AH$r").Value = 0
    $ws.Range("AI$r").Value = 0
    $ws.Range("AJ$r").Value = 0
    $ws.Range("AK$r").Value = 0
    $ws.Range("AL$r").Value = 0
}

# --- sheet view: scroll position and active selection moved --------------
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F7").Select()
